$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.894.32"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.125.78"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.85"
$ws.Range("E5").Value = "  -2.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.99"
$ws.Range("E6").Value = "  -4.95%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.123.20"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.21"
$ws.Range("E11").Value = "  -3.71%  "
$ws.Range("E12").Value = "  -3.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("E13").Value = "  -5.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.90"
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.642.92"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.970.37"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.123.60"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  -4.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.51"
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.04"
$ws.Range("E21").Value = "  -3.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.694"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.87"
$ws.Range("E25").Value = "  -4.25%  "
$ws.Range("E27").Value = "  -1.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.85"
$ws.Range("E28").Value = "  -5.94%  "
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.81"
$ws.Range("E30").Value = "  -5.15%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.51"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("E33").Value = "  -4.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.52"
$ws.Range("E34").Value = "  -4.99%  "
$ws.Range("E35").Value = "  -3.34%  "
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("E37").Value = "  -4.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0676"
$ws.Range("E38").Value = "  -12.04%  "
$ws.Range("E39").Value = "  -2.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "413.38"
$ws.Range("E40").Value = "  -7.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.906.34"
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.14"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.66"
$ws.Range("E43").Value = "  -11.16%  "
$ws.Range("E44").Value = "  -6.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.257"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  -5.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.21"
$ws.Range("E48").Value = "  -2.96%  "
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("E50").Value = "  -8.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.09"
$ws.Range("E51").Value = "  -0.55%  "
